$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.710.49'
$ws.Range("E2").Value = '  +1.40%  '
$ws.Range("D3").Value = '2.104.51'
$ws.Range("E3").Value = '  +5.22%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.82'
$ws.Range("E5").Value = '  +2.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5283'
$ws.Range("E7").Value = '  +3.63%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4348'
$ws.Range("E8").Value = '  +5.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08937'
$ws.Range("E9").Value = '  +2.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '46.66'
$ws.Range("E10").Value = '  +9.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.167'
$ws.Range("E11").Value = '  +2.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.66'
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").Value = '2.106.10'
$ws.Range("E13").Value = '  +5.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.702'
$ws.Range("E14").Value = '  +2.67%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.763'
$ws.Range("E15").Value = '  +4.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '96.94'
$ws.Range("E16").Value = '  +3.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.004'
$ws.Range("E17").Value = '  +0.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001128'
$ws.Range("E18").Value = '  +1.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06669'
$ws.Range("E19").Value = '  +2.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.97'
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.294'
$ws.Range("E22").Value = '  +1.94%  '
$ws.Range("D23").Value = '30.768.16'
$ws.Range("E23").Value = '  +1.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.25'
$ws.Range("E24").Value = '  +3.71%  '
$ws.Range("D25").Value = '2.352.06'
$ws.Range("E25").Value = '  +5.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.289'
$ws.Range("E26").Value = '  +3.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.54'
$ws.Range("E27").Value = '  +0.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.561'
$ws.Range("E28").Value = '  +6.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '161.85'
$ws.Range("E29").Value = '  -0.74%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.71'
$ws.Range("E30").Value = '  +0.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.196'
$ws.Range("E31").Value = '  +4.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1078'
$ws.Range("E32").Value = '  +2.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.139'
$ws.Range("E33").Value = '  +1.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.541'
$ws.Range("E34").Value = '  +15.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.873'
$ws.Range("E35").Value = '  +1.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02590'
$ws.Range("E36").Value = '  +3.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.573'
$ws.Range("E37").Value = '  +6.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.525'
$ws.Range("E38").Value = '  +2.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06728'
$ws.Range("E39").Value = '  +1.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.59'
$ws.Range("E40").Value = '  +2.66%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2269'
$ws.Range("E41").Value = '  +2.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6796'
$ws.Range("E42").Value = '  +2.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.247'
$ws.Range("E43").Value = '  +1.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.002'
$ws.Range("E44").Value = '  +0.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6395'
$ws.Range("E45").Value = '  +3.96%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.95'
$ws.Range("E46").Value = '  +1.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.210'
$ws.Range("E47").Value = '  +0.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.644'
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.254'
$ws.Range("E49").Value = '  -0.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '82.69'
$ws.Range("E50").Value = '  +3.19%  '
$ws.Range("E51").Value = '  +6.93%  '
